$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 133; existing rows 133-156 shift down to 136-159.
$ws.Rows("133:135").Insert()

# Common (fixed) values shared by every row in this product block
$mercadoId  = 1
$mercado    = "Agrícola del Norte S.A. de Arica"
$region     = "Arica y Parinacota"
$codreg     = 15
$tipo       = "Fruta"
$productoId = 100108
$producto   = "Tropicales y subtropicales"
$categoriaId = 100108002
$categoria  = "Mango"

function Set-MangoRow($Row, $Fecha, $Variedad, $Calidad, $Volumen, $PrecioMin, $PrecioMax, $PrecioProm, $Unidad, $Origen, $PrecioKg, $KgUnidad) {
    $ws.Cells.Item($Row, 1).Value = $mercadoId
    $ws.Cells.Item($Row, 2).Value = $mercado
    $ws.Cells.Item($Row, 3).Value = $region
    $ws.Cells.Item($Row, 4).Value = $Fecha
    $ws.Cells.Item($Row, 5).Value = $codreg
    $ws.Cells.Item($Row, 6).Value = $tipo
    $ws.Cells.Item($Row, 7).Value = $productoId
    $ws.Cells.Item($Row, 8).Value = $producto
    $ws.Cells.Item($Row, 9).Value = $categoriaId
    $ws.Cells.Item($Row, 10).Value = $categoria
    $ws.Cells.Item($Row, 11).Value = $Variedad
    $ws.Cells.Item($Row, 12).Value = $Calidad
    $ws.Cells.Item($Row, 13).Value = $Volumen
    $ws.Cells.Item($Row, 14).Value = $PrecioMin
    $ws.Cells.Item($Row, 15).Value = $PrecioMax
    $ws.Cells.Item($Row, 16).Value = $PrecioProm
    $ws.Cells.Item($Row, 17).Value = $Unidad
    $ws.Cells.Item($Row, 18).Value = $Origen
    $ws.Cells.Item($Row, 19).Value = $PrecioKg
    $ws.Cells.Item($Row, 20).Value = $KgUnidad
}

Set-MangoRow 133 44889 "Sin especificar" "Primera" 490 7000 7500 7306 "$/bandeja 4 kilos" "Perú" 1826 4
Set-MangoRow 134 44889 "Sin especificar" "Segunda" 470 7000 7500 7319 "$/bandeja 4 kilos" "Perú" 1830 4
Set-MangoRow 135 44889 "Sin especificar" "Tercera" 400 7000 7500 7375 "$/bandeja 4 kilos" "Perú" 1844 4
